$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet "Tiles" -> "FuelTanks"
$ws.Name = "FuelTanks"

# Update values (replace formula in A2 with a plain literal value)
$ws.Range("A1").Value = 10
$ws.Range("B1").Value = 2
$ws.Range("A2").Value = 20
$ws.Range("B2").Value = 2

# Update selection to A3
$ws.Range("A3").Select()
